# Adds a "2022" column (K) to the 9.5.1 worksheet, mirroring the existing
# "2021" column (J) formatting, and refreshes the sheet's active selection.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 3: blank thick-bottom-border spacer cell, same style as J3 ---
$ws.Range("J3").Copy()
$ws.Range("K3").PasteSpecial(-4122)   # xlPasteFormats

# --- Row 4: year header "2022", same style as J4 ("2021") ---
$ws.Range("J4").Copy()
$ws.Range("K4").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(4, 11).Value = 2022

# --- Row 5: new data point for 2022, same style as J5 ---
$ws.Range("J5").Copy()
$ws.Range("K5").PasteSpecial(-4122)   # xlPasteFormats
$ws.Cells.Item(5, 11).Value = 0.11705180708279034

# --- Column widths: D:K all become a uniform 9-character-wide column
#     (replacing J's old "bestFit" width), matching width="9" in the XML.
$ws.Range("D1:K1").EntireColumn.ColumnWidth = 8.166666666666666

# --- Move the active selection to J12 (was H6) ---
$ws.Range("J12").Select()
